$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Part Of Test Plan") values are flipped: every row becomes "Y"
# except row 20, which becomes "N".
for ($r = 2; $r -le 22; $r++) {
    if ($r -eq 20) {
        $ws.Cells.Item($r, 4).Value = "N"
    } else {
        $ws.Cells.Item($r, 4).Value = "Y"
    }
}
